$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.480.33'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.70'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.73'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3699'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.83'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3387'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.142'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07544'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.27'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.029'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.955'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.572.94'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001120'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06755'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.302'
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.42'
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.22'
$ws.Range('E23').Value = '  +2.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.484.41'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.369'
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.604'
$ws.Range('E26').Value = '  -3.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.05'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.09'
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.21'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.749.11'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.072'
$ws.Range('E32').Value = '  +7.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.254'
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.016'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.762'
$ws.Range('E35').Value = '  -3.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08364'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02489'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.356'
$ws.Range('E38').Value = '  -4.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2299'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06549'
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.441'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6226'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.09'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5860'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.26'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.073'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07332'
$ws.Range('E51').Value = '  +0.12%  '
